$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (172, 173) for date 2025-11-25 (serial 45986),
# one for each station, following the existing per-hour layout in C:Z.
$ws.Cells.Item(172, 1).Value = 45986
$ws.Cells.Item(172, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(172, 3).Value = 562.94899999999996
$ws.Cells.Item(172, 4).Value = 931.83399999999995
$ws.Cells.Item(172, 5).Value = 363.98599999999999
$ws.Cells.Item(172, 6).Value = 461.97899999999998
$ws.Cells.Item(172, 7).Value = 198.97400000000002
$ws.Cells.Item(172, 8).Value = 579.83699999999999
$ws.Cells.Item(172, 9).Value = 506.33500000000009
$ws.Cells.Item(172, 10).Value = 130.81900000000002
$ws.Cells.Item(172, 11).Value = 88.313999999999993
$ws.Cells.Item(172, 12).Value = 110.264
$ws.Cells.Item(172, 13).Value = 180.63
$ws.Cells.Item(172, 14).Value = 195.42199999999997
$ws.Cells.Item(172, 15).Value = 819.18100000000015
$ws.Cells.Item(172, 16).Value = 1255.0550000000003
$ws.Cells.Item(172, 17).Value = 403.10700000000008
$ws.Cells.Item(172, 18).Value = 341.45
$ws.Cells.Item(172, 19).Value = 261.48500000000001
$ws.Cells.Item(172, 20).Value = 200.15600000000001
$ws.Cells.Item(172, 21).Value = 65.84
$ws.Cells.Item(172, 22).Value = 136.16999999999999
$ws.Cells.Item(172, 23).Value = 40.433
$ws.Cells.Item(172, 24).Value = 187.36899999999997
$ws.Cells.Item(172, 25).Value = 40.54
$ws.Cells.Item(172, 26).Value = 82.444000000000003

$ws.Cells.Item(173, 1).Value = 45986
$ws.Cells.Item(173, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(173, 3).Value = 236.98599999999999
$ws.Cells.Item(173, 4).Value = 431.16399999999999
$ws.Cells.Item(173, 5).Value = 44
$ws.Cells.Item(173, 6).Value = 153.95600000000002
$ws.Cells.Item(173, 7).Value = 9.8569999999999993
$ws.Cells.Item(173, 8).Value = 282.56700000000001
$ws.Cells.Item(173, 9).Value = 77.078000000000003
$ws.Cells.Item(173, 10).Value = 130.6
$ws.Cells.Item(173, 11).Value = 224.047
$ws.Cells.Item(173, 12).Value = 304.56800000000004
$ws.Cells.Item(173, 13).Value = 150.26599999999999
$ws.Cells.Item(173, 14).Value = 330.762
$ws.Cells.Item(173, 15).Value = 420.09700000000004
$ws.Cells.Item(173, 16).Value = 340.88099999999997
$ws.Cells.Item(173, 17).Value = 346.60099999999994
$ws.Cells.Item(173, 18).Value = 125.949
$ws.Cells.Item(173, 19).Value = 77.671999999999997
$ws.Cells.Item(173, 20).Value = 78.337999999999994
$ws.Cells.Item(173, 21).Value = 14.867000000000001
$ws.Cells.Item(173, 22).Value = 39.414000000000001
$ws.Cells.Item(173, 23).Value = 79.467999999999989
$ws.Cells.Item(173, 24).Value = 9.8439999999999994
$ws.Cells.Item(173, 25).Value = 90.647000000000006
$ws.Cells.Item(173, 26).Value = 37.515999999999998

# Move the selection/active cell to B7 (matches the saved view state in the diff).
$ws.Range("B7").Select()
